# Auto-generated edit script: update FFXIV Leve profit/price columns (H-N)
# across multiple crafting-class sheets to refreshed market-board figures.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (22 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1110.7778
$ws.Range("I103").Value = 866.6667
$ws.Range("K103").Value = 2600.0001
$ws.Range("M103").Value = -2014.0001
$ws.Range("H129").Value = 975.4681
$ws.Range("J129").Value = 989.2857
$ws.Range("L129").Value = 2967.8571
$ws.Range("N129").Value = -12967.8571
$ws.Range("H132").Value = 1295.4
$ws.Range("I132").Value = 1012.3684
$ws.Range("J132").Value = 2191.6667
$ws.Range("K132").Value = 3037.1052
$ws.Range("L132").Value = 6575.000100000001
$ws.Range("M132").Value = -507.1052
$ws.Range("N132").Value = -11635.0001
$ws.Range("H138").Value = 2287.0132
$ws.Range("I138").Value = 916.2558
$ws.Range("J138").Value = 4073.1516
$ws.Range("K138").Value = 2748.7674
$ws.Range("L138").Value = 12219.4548
$ws.Range("M138").Value = 2391.2326
$ws.Range("N138").Value = -22499.4548

# --- Sheet: ARM (47 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4282.66
$ws.Range("I32").Value = 3037.975
$ws.Range("J32").Value = 9261.4
$ws.Range("K32").Value = 3037.975
$ws.Range("L32").Value = 9261.4
$ws.Range("M32").Value = -2750.975
$ws.Range("N32").Value = -9835.4
$ws.Range("H74").Value = 1107.5435
$ws.Range("I74").Value = 1032.1724
$ws.Range("J74").Value = 1236.1177
$ws.Range("K74").Value = 1032.1724
$ws.Range("L74").Value = 1236.1177
$ws.Range("M74").Value = -158.1723999999999
$ws.Range("N74").Value = -2984.1177
$ws.Range("H77").Value = 1107.5435
$ws.Range("I77").Value = 1032.1724
$ws.Range("J77").Value = 1236.1177
$ws.Range("K77").Value = 5160.861999999999
$ws.Range("L77").Value = 6180.5885
$ws.Range("M77").Value = -792.8619999999992
$ws.Range("N77").Value = -14916.5885
$ws.Range("H97").Value = 1165.6111
$ws.Range("I97").Value = 1004.53845
$ws.Range("J97").Value = 1584.4
$ws.Range("K97").Value = 1004.53845
$ws.Range("L97").Value = 1584.4
$ws.Range("M97").Value = -508.53845
$ws.Range("N97").Value = -2576.4
$ws.Range("H102").Value = 12347279
$ws.Range("I102").Value = 12347279
$ws.Range("K102").Value = 12347279
$ws.Range("M102").Value = -12345657
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H132").Value = 1925601.6
$ws.Range("I132").Value = 1718.2972
$ws.Range("J132").Value = 6671180.5
$ws.Range("K132").Value = 5154.8916
$ws.Range("L132").Value = 20013541.5
$ws.Range("M132").Value = -2624.8916
$ws.Range("N132").Value = -20018601.5
$ws.Range("H133").Value = 38337
$ws.Range("J133").Value = 38337
$ws.Range("L133").Value = 38337
$ws.Range("N133").Value = -43397

# --- Sheet: BSM (32 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1518.375
$ws.Range("I86").Value = 1505.75
$ws.Range("J86").Value = 1543.625
$ws.Range("K86").Value = 1505.75
$ws.Range("L86").Value = 1543.625
$ws.Range("M86").Value = -382.75
$ws.Range("N86").Value = -3789.625
$ws.Range("H89").Value = 1518.375
$ws.Range("I89").Value = 1505.75
$ws.Range("J89").Value = 1543.625
$ws.Range("K89").Value = 7528.75
$ws.Range("L89").Value = 7718.125
$ws.Range("M89").Value = -1912.75
$ws.Range("N89").Value = -18950.125
$ws.Range("H105").Value = 15633133
$ws.Range("I105").Value = 33348178
$ws.Range("J105").Value = 2210
$ws.Range("K105").Value = 33348178
$ws.Range("L105").Value = 2210
$ws.Range("M105").Value = -33346431
$ws.Range("N105").Value = -5704
$ws.Range("H107").Value = 90910500
$ws.Range("I107").Value = 250001470
$ws.Range("K107").Value = 250001470
$ws.Range("M107").Value = -249999550
$ws.Range("H134").Value = 4397.2095
$ws.Range("I134").Value = 6019.9565
$ws.Range("J134").Value = 2531.05
$ws.Range("K134").Value = 18059.8695
$ws.Range("L134").Value = 7593.150000000001
$ws.Range("M134").Value = -15524.8695
$ws.Range("N134").Value = -12663.15

# --- Sheet: CRP (28 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 266804.38
$ws.Range("I31").Value = 1633.973
$ws.Range("J31").Value = 712772.8
$ws.Range("K31").Value = 1633.973
$ws.Range("L31").Value = 712772.8
$ws.Range("M31").Value = -1338.973
$ws.Range("N31").Value = -713362.8
$ws.Range("H34").Value = 266804.38
$ws.Range("I34").Value = 1633.973
$ws.Range("J34").Value = 712772.8
$ws.Range("K34").Value = 1633.973
$ws.Range("L34").Value = 712772.8
$ws.Range("M34").Value = -1431.973
$ws.Range("N34").Value = -713176.8
$ws.Range("H132").Value = 2086.93
$ws.Range("I132").Value = 1746.5385
$ws.Range("J132").Value = 2824.4443
$ws.Range("K132").Value = 5239.6155
$ws.Range("L132").Value = 8473.332900000001
$ws.Range("M132").Value = -2709.6155
$ws.Range("N132").Value = -13533.3329
$ws.Range("H134").Value = 1489.9556
$ws.Range("I134").Value = 1612.1333
$ws.Range("J134").Value = 1245.6
$ws.Range("K134").Value = 4836.3999
$ws.Range("L134").Value = 3736.8
$ws.Range("M134").Value = -2301.3999
$ws.Range("N134").Value = -8806.799999999999

# --- Sheet: CUL (15 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 119.8
$ws.Range("I11").Value = 119.8
$ws.Range("K11").Value = 359.4
$ws.Range("M11").Value = -219.4
$ws.Range("H129").Value = 37038720
$ws.Range("J129").Value = 2758.25
$ws.Range("L129").Value = 8274.75
$ws.Range("N129").Value = -18274.75
$ws.Range("H131").Value = 2942086
$ws.Range("I131").Value = 12500425
$ws.Range("J131").Value = 1058.5
$ws.Range("K131").Value = 37501275
$ws.Range("L131").Value = 3175.5
$ws.Range("M131").Value = -37496235
$ws.Range("N131").Value = -13255.5

# --- Sheet: GSM (29 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2963.6365
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 3950
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 3950
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -5946
$ws.Range("H83").Value = 2963.6365
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 3950
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 19750
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -29734
$ws.Range("H97").Value = 1003.75
$ws.Range("I97").Value = 1003.75
$ws.Range("K97").Value = 1003.75
$ws.Range("M97").Value = -507.75
$ws.Range("H102").Value = 1827.5834
$ws.Range("I102").Value = 1808
$ws.Range("J102").Value = 1834.1111
$ws.Range("K102").Value = 1808
$ws.Range("L102").Value = 1834.1111
$ws.Range("M102").Value = -186
$ws.Range("N102").Value = -5078.1111
$ws.Range("H122").Value = 40956308
$ws.Range("I122").Value = 66552470
$ws.Range("K122").Value = 199657410
$ws.Range("M122").Value = -199654960

# --- Sheet: LTW (20 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 145057.72
$ws.Range("I7").Value = 168484
$ws.Range("K7").Value = 168484
$ws.Range("M7").Value = -168372
$ws.Range("H108").Value = 322626
$ws.Range("J108").Value = 322626
$ws.Range("L108").Value = 322626
$ws.Range("N108").Value = -330306
$ws.Range("H126").Value = 145057.72
$ws.Range("I126").Value = 168484
$ws.Range("K126").Value = 505452
$ws.Range("M126").Value = -502982
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 6686.648
$ws.Range("I136").Value = 4106
$ws.Range("K136").Value = 12318
$ws.Range("M136").Value = -9768

# --- Sheet: WVR (11 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 80001300
$ws.Range("J107").Value = 5714808
$ws.Range("L107").Value = 17144424
$ws.Range("N107").Value = -17148264
$ws.Range("H132").Value = 918.38464
$ws.Range("I132").Value = 650.5
$ws.Range("J132").Value = 1989.9231
$ws.Range("K132").Value = 1951.5
$ws.Range("L132").Value = 5969.7693
$ws.Range("M132").Value = 578.5
$ws.Range("N132").Value = -11029.7693

